$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "2020" header style (J4/K4 should look like I4 but without
# the fill flag / table look - we fix that up below) and the corresponding
# data-column styles used for rows 5-14, then fill in 2021 / 2022 data.

# --- Header row 4: years 2021, 2022 ---
$ws.Range("I4").Copy() | Out-Null
$ws.Range("J4:K4").PasteSpecial(-4104) | Out-Null
$ws.Range("J4").Value = 2021
$ws.Range("K4").Value = 2022
$ws.Range("J4:K4").Borders.Item(10).Weight = 2
$ws.Range("J4:K4").WrapText = $true
$ws.Range("J4:K4").HorizontalAlignment = -4152
$ws.Range("J4:K4").VerticalAlignment = -4160

# --- Row 5 (first data row, percent-like 165 number format) ---
$ws.Range("I5").Copy() | Out-Null
$ws.Range("J5:K5").PasteSpecial(-4104) | Out-Null
$ws.Range("J5").Value = 34.075233127500141
$ws.Range("K5").Value = 35.305353068702679

# --- Rows 6-13 (regional rows) ---
$dataRows = @{
    6  = @(44.487602536118636, 49.31549563692068)
    7  = @(40.668697007891453, 45.444207273635158)
    8  = @(50.797011639929529, 46.810603774236895)
    9  = @(46.848562449074493, 45.450816127137941)
    10 = @(44.458036086558309, 50.474514452886076)
    11 = @(40.532201616746903, 40.14796186663478)
    12 = @(33.353175884696697, 29.735683954543184)
    13 = @(10.46405303463253, 12.912087912087852)
}

foreach ($r in $dataRows.Keys) {
    $ws.Range("I$r").Copy() | Out-Null
    $ws.Range("J${r}:K${r}").PasteSpecial(-4104) | Out-Null
    $vals = $dataRows[$r]
    $ws.Range("J$r").Value = $vals[0]
    $ws.Range("K$r").Value = $vals[1]
}

# --- Row 14 (bottom total row with thicker bottom border) ---
$ws.Range("I14").Copy() | Out-Null
$ws.Range("J14:K14").PasteSpecial(-4104) | Out-Null
$ws.Range("J14").Value = 43.479082661290349
$ws.Range("K14").Value = 41.117034465658314

$excel.CutCopyMode = 0

# --- Sheet-level bookkeeping to mirror the diff ---
$ws.Range("M6").Select()
